$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The duplicate survey question column (K) is removed; everything to its
# right (L:R) shifts one column to the left (becomes K:Q).
$ws.Columns("K").Delete()

# Restore/update the view state recorded for the sheet after the edit.
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("K5").Select()
